# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-29 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 3
    4  = 3
    5  = 8
    6  = 3
    7  = 6
    8  = 6
    9  = 4
    10 = 6
    11 = 3
    12 = 8
    13 = 5
    14 = 4
    15 = 6
    16 = 5
    17 = 5
    18 = 7
    19 = 7
    20 = 7
    21 = 5
    22 = 5
    23 = 5
    24 = 3
    25 = 2
    26 = 3
    27 = 5
    28 = 1
    29 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
